$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '45.075.71'
$ws.Cells.Item(2, 5).Value = '  +3.86%  '
$ws.Cells.Item(3, 4).Value = '2.428.61'
$ws.Cells.Item(3, 5).Value = '  +0.79%  '
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '317.00'
$ws.Cells.Item(5, 5).Value = '  +3.68%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '102.65'
$ws.Cells.Item(6, 5).Value = '  +5.58%  '
$ws.Cells.Item(7, 5).Value = '  +1.24%  '
$ws.Cells.Item(8, 5).Value = '  -0.10%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.526'
$ws.Cells.Item(9, 5).Value = '  +7.61%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '35.41'
$ws.Cells.Item(10, 5).Value = '  +1.30%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0801'
$ws.Cells.Item(11, 5).Value = '  +0.97%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.122'
$ws.Cells.Item(12, 5).Value = '  -2.42%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '18.14'
$ws.Cells.Item(13, 5).Value = '  -1.87%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '7.02'
$ws.Cells.Item(14, 5).Value = '  +1.99%  '
$ws.Cells.Item(15, 4).Value = '2.809.99'
$ws.Cells.Item(15, 5).Value = '  +0.94%  '
$ws.Cells.Item(16, 4).Value = '2.425.54'
$ws.Cells.Item(16, 5).Value = '  +0.10%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.837'
$ws.Cells.Item(17, 5).Value = '  +1.60%  '
$ws.Cells.Item(18, 4).Value = '44.997.23'
$ws.Cells.Item(18, 5).Value = '  +3.75%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '12.25'
$ws.Cells.Item(19, 5).Value = '  +0.92%  '
$ws.Cells.Item(20, 5).Value = '  -0.71%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0922'
$ws.Cells.Item(21, 5).Value = '  +2.46%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '68.81'
$ws.Cells.Item(22, 5).Value = '  +0.47%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '243.91'
$ws.Cells.Item(23, 5).Value = '  +2.59%  '
$ws.Cells.Item(24, 5).Value = '  +1.36%  '
$ws.Cells.Item(25, 5).Value = '  +1.80%  '
$ws.Cells.Item(26, 5).Value = '  +0.00%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '25.29'
$ws.Cells.Item(27, 5).Value = '  +1.56%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.28'
$ws.Cells.Item(28, 5).Value = '  +3.06%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '9.55'
$ws.Cells.Item(29, 5).Value = '  +1.57%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '49.29'
$ws.Cells.Item(30, 5).Value = '  +3.39%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '32.75'
$ws.Cells.Item(31, 5).Value = '  +1.61%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '20.18'
$ws.Cells.Item(32, 5).Value = '  +10.17%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.125'
$ws.Cells.Item(33, 5).Value = '  +9.91%  '
$ws.Cells.Item(34, 5).Value = '  +1.96%  '
$ws.Cells.Item(35, 5).Value = '  +0.28%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.0762'
$ws.Cells.Item(36, 5).Value = '  +2.44%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.87'
$ws.Cells.Item(37, 5).Value = '  -0.40%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '4.42'
$ws.Cells.Item(38, 5).Value = '  +1.11%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '2.87'
$ws.Cells.Item(39, 5).Value = '  -2.27%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '125.30'
$ws.Cells.Item(40, 5).Value = '  -5.67%  '
$ws.Cells.Item(41, 5).Value = '  -2.44%  '
$ws.Cells.Item(42, 5).Value = '  +0.95%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '20.68'
$ws.Cells.Item(43, 5).Value = '  -3.58%  '
$ws.Cells.Item(44, 5).Value = '  +1.91%  '
$ws.Cells.Item(45, 4).Value = '1.933.86'
$ws.Cells.Item(45, 5).Value = '  -0.76%  '
$ws.Cells.Item(46, 5).Value = '  -2.29%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.92'
$ws.Cells.Item(47, 5).Value = '  +4.16%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '9.24'
$ws.Cells.Item(48, 5).Value = '  -0.26%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.81'
$ws.Cells.Item(49, 5).Value = '  +17.34%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '76.44'
$ws.Cells.Item(50, 5).Value = '  +5.85%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '53.88'
$ws.Cells.Item(51, 5).Value = '  +2.67%  '
